# Apply the diff: update dSF (column F) values for specific rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = -3
$ws.Range("F14").Value = -3
$ws.Range("F16").Value = -2
$ws.Range("F27").Value = -3
$ws.Range("F28").Value = -13
$ws.Range("F31").Value = 4
